$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (AngleSensor): add Developer / finished / Mode / interface info
$ws.Range("D6").Value = "Lawrie"
$ws.Range("E6").Value = "N"
$ws.Range("F6").Value = "Angle"
$ws.Range("G6").Value = "SamplerProvider"

# Row 7 (BarometricHTSensor): add Developer / finished / Mode / interface info
$ws.Range("D7").Value = "Lawrie"
$ws.Range("E7").Value = "N"
$ws.Range("F7").Value = "Pressure"
$ws.Range("G7").Value = "SampleProvider"

# Rename "CalibratedSampleProvider" occurrences to "Calibrate, SampleProvider"
$ws.Range("G10").Value = "Calibrate, SampleProvider"
$ws.Range("G11").Value = "Calibrate, SampleProvider"

# Update the active selection to match the committed view state
$ws.Range("G7").Select()
